$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting from the row above into the new row's date cell
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)

# New race result row (Volta a la Comunitat Valenciana, Stage 1)
$ws.Range("A9").Value = "2/4/2026"
$ws.Range("B9").Value = "Volta a la Comunitat Valenciana"
$ws.Range("C9").Value = "Stage 1"
$ws.Range("D9").Value = "Biniam Girmay"
$ws.Range("E9").Value = "Arne Marit"
$ws.Range("F9").Value = "Giovanni Lonardi"
$ws.Range("G9").Value = "Carl-Frederik Bévort"
$ws.Range("H9").Value = "Aleksandr Vlasov"
$ws.Range("I9").Value = "Alberto Bruttomesso"
$ws.Range("J9").Value = "Tommaso Nencini"
$ws.Range("K9").Value = "Mikel Retegi"
$ws.Range("L9").Value = "Dries Van Gestel"
$ws.Range("M9").Value = "Clément Alleno"

# Printed page orientation was set to portrait
$ws.PageSetup.Orientation = 1

$ws.Range("E25:F25").Select() | Out-Null
